# "Add files via upload" — re-uploaded workbook renames the sole worksheet
# from "Аркуш1" (Ukrainian) to "Лист1" and resets the saved cursor/selection
# back to the default cell (A1) instead of the previously recorded I12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab.
$ws.Name = "Лист1"

# Make sure the sheet is the active one, then park the selection on A1 so the
# view no longer points at the old "I12" cursor position.
$ws.Activate()
$ws.Range("A1").Select()
